$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: copy the number-format/border/fill look of the same-row D cell
#     into the E cell, then stamp the "0.0" number format with right/center
#     alignment on top (mirrors how the 2020 column was filled in next to
#     the existing 2018 column). ---
function Fill-NumericE($row, $value) {
    $d = $ws.Range("D$row")
    $e = $ws.Range("E$row")
    $d.Copy()
    $e.PasteSpecial(-4122)
    $e.NumberFormat = "0.0"
    $e.HorizontalAlignment = -4152
    $e.VerticalAlignment = -4108
    if ($value -ne $null) {
        $e.Value = $value
    }
}

function Fill-DashE($row) {
    $d = $ws.Range("D$row")
    $e = $ws.Range("E$row")
    $d.Copy()
    $e.PasteSpecial(-4122)
    $e.NumberFormat = "0.0"
    $e.HorizontalAlignment = -4152
    $e.VerticalAlignment = -4108
    $e.Value = "-"
}

# Header: 2020 column header, same look as the 2018 header in D3
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = 2020

# Data rows 4-18 (numeric, some blank)
Fill-NumericE 4 22.1
Fill-NumericE 5 $null
Fill-NumericE 6 52.7
Fill-NumericE 7 5
Fill-NumericE 8 $null
Fill-NumericE 9 4.8
Fill-NumericE 10 15.8
Fill-NumericE 11 13.5
Fill-NumericE 12 9.6
Fill-NumericE 13 2.7
Fill-NumericE 14 14.7
Fill-NumericE 15 18.2
Fill-NumericE 16 74
Fill-NumericE 17 35.1
Fill-NumericE 18 $null

# Rows 19-23: education-of-head-of-household block -> not available, dash
Fill-DashE 19
Fill-DashE 20
Fill-DashE 21
Fill-DashE 22
Fill-DashE 23

# Row 24: wealth-quintile header row, stays blank
Fill-NumericE 24 $null

# Rows 25-29: wealth quintile block -> not available, dash
Fill-DashE 25
Fill-DashE 26
Fill-DashE 27
Fill-DashE 28
Fill-DashE 29

# Restore the active selection the workbook was left on
$ws.Range("J24").Select()
